$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.362.92"
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("D3").Value = "2.976.04"
$ws.Range("E3").Value = "  +2.28%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "380.11"
$ws.Range("E5").Value = "  +3.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.47"
$ws.Range("E6").Value = "  +2.21%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.594"
$ws.Range("E9").Value = "  +0.39%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.32"
$ws.Range("E10").Value = "  +1.19%  "
$ws.Range("E11").Value = "  +0.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0844"
$ws.Range("E12").Value = "  +0.97%  "
$ws.Range("D13").Value = "3.436.42"
$ws.Range("E13").Value = "  +2.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.43"
$ws.Range("E14").Value = "  -0.18%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.56"
$ws.Range("E15").Value = "  +2.34%  "
$ws.Range("D16").Value = "2.969.33"
$ws.Range("E16").Value = "  +2.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.966"
$ws.Range("E17").Value = "  +2.44%  "
$ws.Range("D18").Value = "51.305.30"
$ws.Range("E18").Value = "  +0.35%  "
$ws.Range("E19").Value = "  +2.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.40"
$ws.Range("E20").Value = "  +2.31%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.92"
$ws.Range("E21").Value = "  +0.88%  "
$ws.Range("E22").Value = "  +1.74%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.45"
$ws.Range("E23").Value = "  +1.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "261.73"
$ws.Range("E24").Value = "  +0.67%  "
$ws.Range("E25").Value = "  +5.23%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.50"
$ws.Range("E26").Value = "  +23.54%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.70"
$ws.Range("E27").Value = "  +10.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.172"
$ws.Range("E28").Value = "  +0.14%  "
$ws.Range("B29").Value = "Hedera"
$ws.Range("C29").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.112"
$ws.Range("E29").Value = "  +10.27%  "
$ws.Range("B30").Value = "Dai"
$ws.Range("C30").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.92"
$ws.Range("E31").Value = "  +0.73%  "
$ws.Range("E32").Value = "  -0.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "35.27"
$ws.Range("E33").Value = "  +1.46%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "51.28"
$ws.Range("E34").Value = "  +0.63%  "
$ws.Range("E35").Value = "  -2.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0446"
$ws.Range("E36").Value = "  +6.45%  "
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("E38").Value = "  +0.40%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "17.24"
$ws.Range("E39").Value = "  +0.82%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.60"
$ws.Range("E40").Value = "  -1.67%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.85"
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("E42").Value = "  +2.44%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "124.64"
$ws.Range("E43").Value = "  +4.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.98"
$ws.Range("E44").Value = "  -0.43%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.291"
$ws.Range("E45").Value = "  +22.95%  "
$ws.Range("E46").Value = "  -0.98%  "
$ws.Range("E47").Value = "  +3.41%  "
$ws.Range("D48").Value = "2.044.31"
$ws.Range("E48").Value = "  +1.23%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.25"
$ws.Range("E49").Value = "  +2.40%  "
$ws.Range("E50").Value = "  +10.53%  "
$ws.Range("E51").Value = "  +1.56%  "
